# Swap the full contents of row 2 and row 3 on the active sheet.
# (The two observation records were swapped in the source export: the
# row that used to be "A 103410817 / Flakaträsk ..." becomes the
# "A 7125672 / Vänjaurbäck ..." record, and vice versa.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 51   # column AY

# Columns whose text values could otherwise be auto-coerced by Excel's
# smart entry (numeric-looking "1", date-looking "2022-08-08", ...) and
# therefore need to be forced to stay text while we round-trip them.
$textForceCols = @(9, 25, 27)   # I, Y, AA

# --- snapshot both rows ------------------------------------------------
$row2Vals = @()
$row3Vals = @()

for ($c = 1; $c -le $lastCol; $c++) {
    $row2Vals += , $ws.Cells.Item(2, $c).Value2
    $row3Vals += , $ws.Cells.Item(3, $c).Value2
}

# --- write row 2 <- old row 3, row 3 <- old row 2 ----------------------
for ($c = 1; $c -le $lastCol; $c++) {
    $cell2 = $ws.Cells.Item(2, $c)
    $cell3 = $ws.Cells.Item(3, $c)

    if ($textForceCols -contains $c) {
        $cell2.NumberFormat = "@"
        $cell3.NumberFormat = "@"
    }

    $cell2.Value2 = $row3Vals[$c - 1]
    $cell3.Value2 = $row2Vals[$c - 1]

    if ($textForceCols -contains $c) {
        $cell2.Style = "Normal"
        $cell3.Style = "Normal"
    }
}
